$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '315.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.21%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '39.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.51%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.130'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.13%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08168'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.67%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.978'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.73%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.378'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '3.84%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '8.341'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.48%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9376'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.96%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1296'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-7.91%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1962'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '1.97%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09027'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-2.04%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03498'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.52%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09760'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.48%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001423'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '2.24%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005959'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.71%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.649'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-6.17%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.319'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.90%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3490'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.83%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1318'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '0.43%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.988'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '7.42%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2493'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.26%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04363'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.65%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001245'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '2.54%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004765'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '9.44%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.04%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-7.62%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02209'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '8.44%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05179'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2.64%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007752'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '4.98%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01041'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '5.41%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1402'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2.71%'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-4.13%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.009292'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.31%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006934'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '9.13%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.08%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002885'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '0.13%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '30.17%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002104'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.08%'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.08%'
